$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:B9 with new computed values
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0.2781908756
$ws.Range("B4").Value = 0.3703939161
$ws.Range("B5").Value = 0.6083302473
$ws.Range("B6").Value = 0.8345726268
$ws.Range("B7").Value = 1.185616626
$ws.Range("B8").Value = 1.742154268
$ws.Range("B9").Value = 0.6689214102

# Clear contents of A10:B15 (keep formatting/style)
$ws.Range("A10:B15").ClearContents()

# Update selection
$ws.Range("A2:A9").Select()
